# Updated symbol list (crypto price/volume refresh).
# Price cells in column D are stored as text, so numeric-looking
# replacements are prefixed with a leading apostrophe to force Excel
# to keep them as text instead of auto-converting to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.73"
$ws.Range("D3").Value = "'21.77"
$ws.Range("D4").Value = "'5.437"
$ws.Range("D5").Value = "'0.05681"
$ws.Range("D6").Value = "'3.378"
$ws.Range("D7").Value = "'0.8073"
$ws.Range("D8").Value = "'1.021"
$ws.Range("D9").Value = "'0.1460"
$ws.Range("D10").Value = "'0.07685"
$ws.Range("D11").Value = "'0.03163"
$ws.Range("D12").Value = "'0.03029"
$ws.Range("D13").Value = "'0.09259"
$ws.Range("D14").Value = "'3.533"
$ws.Range("D15").Value = "'0.001626"
$ws.Range("D16").Value = "'0.04702"
$ws.Range("D17").Value = "'0.01163"
$ws.Range("E17").Value = "16OneONEBestin24h"
$ws.Range("D18").Value = "'0.006353"
$ws.Range("D19").Value = "'0.005025"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("D20").Value = "'0.001044"
$ws.Range("D22").Value = "'0.0003103"
$ws.Range("D24").Value = "'6.425"
$ws.Range("D25").Value = "'2.166"
$ws.Range("D26").Value = "'0.3293"
$ws.Range("D40").Value = "'0.04063"
$ws.Range("D41").Value = "'0.006956"
$ws.Range("D42").Value = "'0.1040"
$ws.Range("D43").Value = "'0.003147"
$ws.Range("D44").Value = "'0.008530"
$ws.Range("D45").Value = "'0.00005930"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.0005507"
$ws.Range("D48").Value = "'0.6831"
$ws.Range("D49").Value = "'0.007982"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D51").Value = "'0.01011"
